# Applies the "Added login page tests" edit:
#  - Renames the TestMethod values in column B (rows 2-5) so that they
#    correctly correspond to the Username/Password validity combination
#    described by columns C and D, and prefixes them with "...Feature...".
#  - Re-selects cell B3 (matching the saved selection state in the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "validateLoginFeatureWithValidUserNameAndValidPassword"
$ws.Range("B3").Value = "validateLoginFeatureWIthInvalidUsernameAndValidPassword"
$ws.Range("B4").Value = "validateLoginFeatureWIthValidUsernameAndInvalidPassword"
$ws.Range("B5").Value = "validateLoginFeatureWIthInvalidUsernameAndInvalidPassword"

$ws.Range("B3").Select()
